$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Class 2 (rows 3-12) ---
# target bandwidth raised from 100 kbps to 600000 kbps
$ws.Range("K3").Value = 600000

# K5 used to mirror C7 (iperf packet size 106 bytes); now computed directly
# from the 1542-byte "max ether packet" constant, same as C10/K16 below.
$ws.Range("K5").Formula = "=1542/1000"

# maxFrameSize (C10) switches from the 1514-byte iperf packet to the
# 1542-byte max-ether-packet constant.
$ws.Range("C10").Formula = "=1542/1000"

# --- Class 1 (rows 14-22) ---
$ws.Range("K14").Value = 500000

# K16 used to mirror C8 (1514/1000); now uses the 1542/1000 constant directly.
$ws.Range("K16").Formula = "=1542/1000"

# --- Class 0 / Best effort (rows 25-33) ---
$ws.Range("K25").Value = 999999
$ws.Range("O25").Value = 1000000

# --- formatting: B34:B37 (the tc qdisc command strings) get an explicit
# numeric-looking "0" number format (was General) ---
$ws.Range("B34:B37").NumberFormat = "0"

# row 37 no longer needs the taller wrapped height now that the rendered
# command text fits back on the default row height
$ws.Rows.Item(37).RowHeight = 12.8

# --- view: scroll/selection moved from O33 to I23 ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("I23").Select() | Out-Null
